$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:O2").Value = 84.375

$ws.Range("P2").Value = 76.82291666666666
$ws.Range("Q2").Value = 75.78125
$ws.Range("R2").Value = 75.78125
$ws.Range("S2").Value = 75.52083333333334
$ws.Range("T2").Value = 75.26041666666666
$ws.Range("U2").Value = 75.26041666666666
